$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 482, shifting existing rows 482-512 down to 483-513.
$ws.Rows.Item(482).Insert()

# Populate the new row 482 with the new record's data.
$ws.Cells.Item(482, 1).Value = 6
$ws.Cells.Item(482, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(482, 3).Value = "Metropolitana"
$ws.Cells.Item(482, 4).Value = 44714
$ws.Cells.Item(482, 5).Value = 13
$ws.Cells.Item(482, 6).Value = 100112044
$ws.Cells.Item(482, 7).Value = "Perejil"
$ws.Cells.Item(482, 8).Value = "Sin especificar"
$ws.Cells.Item(482, 9).Value = "Primera"
$ws.Cells.Item(482, 10).Value = 280
$ws.Cells.Item(482, 11).Value = 8500
$ws.Cells.Item(482, 12).Value = 9000
$ws.Cells.Item(482, 13).Value = 8732
$ws.Cells.Item(482, 14).Value = "$/docena de atados"
$ws.Cells.Item(482, 15).Value = "Región Metropolitana"
$ws.Cells.Item(482, 16).Value = 2911
$ws.Cells.Item(482, 17).Value = 3
$ws.Cells.Item(482, 18).Value = "Hortaliza"
